$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("M26").Value = 20000
$ws.Range("N26").Value = -20688

$ws.Range("H44").Value = 24285.715
$ws.Range("J44").Value = 24285.715
$ws.Range("L44").Value = 24285.715
$ws.Range("N44").Value = -25209.715

$ws.Range("H135").Value = 1028.5294
$ws.Range("I135").Value = 933.46155
$ws.Range("J135").Value = 1337.5
$ws.Range("K135").Value = 8401.15395
$ws.Range("L135").Value = 12037.5
$ws.Range("M135").Value = -5866.15395
$ws.Range("N135").Value = -17107.5

$ws.Range("H137").Value = 53406.137
$ws.Range("I137").Value = 1886.5555
$ws.Range("J137").Value = 76589.95
$ws.Range("K137").Value = 5659.666499999999
$ws.Range("L137").Value = 229769.85
$ws.Range("M137").Value = -3109.666499999999
$ws.Range("N137").Value = -234869.85

$ws.Range("H138").Value = 2200.8
$ws.Range("I138").Value = 2200.8
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 6602.400000000001
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -1462.400000000001

$ws.Range("H141").Value = 39239.582
$ws.Range("I141").Value = 39239.582
$ws.Range("K141").Value = 117718.746
$ws.Range("M141").Value = -112538.746

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15308.833
$ws.Range("I45").Value = 15541.333
$ws.Range("K45").Value = 15541.333
$ws.Range("M45").Value = -15164.333

$ws.Range("H132").Value = 19812.467
$ws.Range("I132").Value = 22497.75
$ws.Range("J132").Value = 9071.333000000001
$ws.Range("K132").Value = 67493.25
$ws.Range("L132").Value = 27213.999
$ws.Range("M132").Value = -64963.25
$ws.Range("N132").Value = -32273.999

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7189.7896
$ws.Range("I134").Value = 6788.625
$ws.Range("J134").Value = 9329.333000000001
$ws.Range("K134").Value = 20365.875
$ws.Range("L134").Value = 27987.999
$ws.Range("M134").Value = -17830.875
$ws.Range("N134").Value = -33057.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1380
$ws.Range("I19").Value = 893.3333
$ws.Range("J19").Value = 1866.6666
$ws.Range("K19").Value = 893.3333
$ws.Range("L19").Value = 1866.6666
$ws.Range("M19").Value = -723.3333
$ws.Range("N19").Value = -2206.6666

$ws.Range("H24").Value = 1380
$ws.Range("I24").Value = 893.3333
$ws.Range("J24").Value = 1866.6666
$ws.Range("K24").Value = 893.3333
$ws.Range("L24").Value = 1866.6666
$ws.Range("M24").Value = -723.3333
$ws.Range("N24").Value = -2206.6666

$ws.Range("H31").Value = 3035.016
$ws.Range("I31").Value = 2248.3872
$ws.Range("J31").Value = 3821.6453
$ws.Range("K31").Value = 2248.3872
$ws.Range("L31").Value = 3821.6453
$ws.Range("M31").Value = -1953.3872
$ws.Range("N31").Value = -4411.6453

$ws.Range("H34").Value = 3035.016
$ws.Range("I34").Value = 2248.3872
$ws.Range("J34").Value = 3821.6453
$ws.Range("K34").Value = 2248.3872
$ws.Range("L34").Value = 3821.6453
$ws.Range("M34").Value = -2046.3872
$ws.Range("N34").Value = -4225.6453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 748.5
$ws.Range("I5").Value = 611.8570999999999
$ws.Range("J5").Value = 885.1429000000001
$ws.Range("K5").Value = 1835.5713
$ws.Range("L5").Value = 2655.4287
$ws.Range("M5").Value = -1723.5713
$ws.Range("N5").Value = -2879.4287

$ws.Range("H22").Value = 2048.75
$ws.Range("I22").Value = 299.5
$ws.Range("J22").Value = 2398.6
$ws.Range("K22").Value = 898.5
$ws.Range("L22").Value = 7195.799999999999
$ws.Range("M22").Value = -729.5
$ws.Range("N22").Value = -7533.799999999999

$ws.Range("H27").Value = 2048.75
$ws.Range("I27").Value = 299.5
$ws.Range("J27").Value = 2398.6
$ws.Range("K27").Value = 898.5
$ws.Range("L27").Value = 7195.799999999999
$ws.Range("M27").Value = -796.5
$ws.Range("N27").Value = -7399.799999999999

$ws.Range("H68").Value = 2505929.5
$ws.Range("J68").Value = 3853044.8
$ws.Range("L68").Value = 11559134.4
$ws.Range("N68").Value = -11560756.4

$ws.Range("H71").Value = 2505929.5
$ws.Range("J71").Value = 3853044.8
$ws.Range("L71").Value = 34677403.2
$ws.Range("N71").Value = -34685515.2

$ws.Range("H92").Value = 839.2857
$ws.Range("I92").Value = 497.5
$ws.Range("J92").Value = 976
$ws.Range("K92").Value = 1492.5
$ws.Range("L92").Value = 2928
$ws.Range("M92").Value = -244.5
$ws.Range("N92").Value = -5424

$ws.Range("H107").Value = 2309.1428
$ws.Range("I107").Value = 1672.6666
$ws.Range("J107").Value = 2786.5
$ws.Range("K107").Value = 5017.9998
$ws.Range("L107").Value = 8359.5
$ws.Range("M107").Value = -3097.9998
$ws.Range("N107").Value = -12199.5

$ws.Range("H135").Value = 748.5
$ws.Range("I135").Value = 611.8570999999999
$ws.Range("J135").Value = 885.1429000000001
$ws.Range("K135").Value = 5506.7139
$ws.Range("L135").Value = 7966.2861
$ws.Range("M135").Value = -2971.7139
$ws.Range("N135").Value = -13036.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

$ws.Range("H132").Value = 5534.7896
$ws.Range("I132").Value = 5619
$ws.Range("J132").Value = 5299
$ws.Range("K132").Value = 16857
$ws.Range("L132").Value = 15897
$ws.Range("M132").Value = -14327
$ws.Range("N132").Value = -20957

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11390.105
$ws.Range("I132").Value = 14454.846
$ws.Range("J132").Value = 4749.8335
$ws.Range("K132").Value = 43364.538
$ws.Range("L132").Value = 14249.5005
$ws.Range("M132").Value = -40834.538
$ws.Range("N132").Value = -19309.5005

$ws.Range("H136").Value = 5344.1177
$ws.Range("I136").Value = 5606.6665
$ws.Range("J136").Value = 3375
$ws.Range("K136").Value = 16819.9995
$ws.Range("L136").Value = 10125
$ws.Range("M136").Value = -14269.9995
$ws.Range("N136").Value = -15225

